$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.270.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.811.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.810.35'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('E13').Value = '  -5.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.453.72'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.813.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.284.66'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '471.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.62%  '
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.960.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.38%  '
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.48'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.36'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.766.87'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.96%  '
$ws.Range('E40').Value = '  +1.37%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '44.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +11.52%  '
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '402.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.88%  '
